$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 18
$ws.Cells.Item($row, 1).Value = 42622.885798611111
$ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item($row, 2).Value = 38
$ws.Cells.Item($row, 3).Value = 64
$ws.Cells.Item($row, 4).Value = 33
$ws.Cells.Item($row, 5).Value = 64
$ws.Cells.Item($row, 6).Value = 20
$ws.Cells.Item($row, 7).Value = 20324
$ws.Cells.Item($row, 8).Value = 19714
$ws.Cells.Item($row, 9).Value = 3139
$ws.Cells.Item($row, 10).Value = 463
$ws.Cells.Item($row, 11).Value = 238
$ws.Cells.Item($row, 12).Value = 51
$ws.Cells.Item($row, 13).Value = 13
$ws.Cells.Item($row, 14).Value = "Noun"
